$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Scenarios")

$ws.Range("D3").Value = "FAIL"
$ws.Range("D4").Value = "FAIL"
$ws.Range("D5").Value = "FAIL"
